# Fixing generated Schema test with generated PredefiendSteps anyType
#
# 1) "Types" sheet gains a new "anyTypeField" / "anyType" row (row 11),
#    mirroring the existing field/type rows above it.
# 2) The active sheet/selection moves from "Formula" (last sheet) back to
#    "Types" (first sheet), with the active cell on "Types" set to C12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Types")

$ws.Range("A11").Value = "field"
$ws.Range("B11").Value = "anyTypeField"
$ws.Range("C11").Value = "anyType"

$ws.Activate()
[void]$ws.Range("C12").Select()
